$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questions")

# Insert a new column before column C (shifts old C->D, D->E), preserving data.
$ws.Range("C1").EntireColumn.Insert()

# New header for inserted column C (row 2) and sequential index values (rows 3-40).
$ws.Cells.Item(2, 3).Value = "I_QSTN"
for ($r = 3; $r -le 40; $r++) {
    $ws.Cells.Item($r, 3).Value = $r - 2
}

# Selection / active sheet bookkeeping to match the authored workbook view state.
$ws.Activate()
$ws.Range("C3:C40").Select()

# Page setup (portrait) now present on the Questions sheet.
$ws.PageSetup.Orientation = 1
